$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values for columns B:E
$ws.Range("B2").Value = 1.352655904642404
$ws.Range("C2").Value = 15.770593801060102
$ws.Range("D2").Value = 20.959639185426568
$ws.Range("E2").Value = 28.479259450375366

# Update row 3 (STR) values for columns B:E
$ws.Range("B3").Value = -19.535092476565751
$ws.Range("C3").Value = 13.102229730301303
$ws.Range("D3").Value = 43.009426650320847
$ws.Range("E3").Value = 19.578343119659678

# Update the selection to match the new selected range B1:E3
$ws.Range("B1:E3").Select() | Out-Null
